$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting setup -----------------------------------------------------
# Rows 19-20 reuse the "section header / header value" look of rows 8-9
# (A col = orange label style, B col = orange centered value style).
$ws.Range("A8:B9").Copy()
$ws.Range("A19:B20").PasteSpecial(-4122)

# Rows 21-24, column B reuse the centered data style used in column B of
# rows 3-6 (style 5); column A keeps the default (no) style.
$ws.Range("B3:B6").Copy()
$ws.Range("B21:B24").PasteSpecial(-4122)

# --- Values -----------------------------------------------------------------
# [Roles] sub-table: name[Role*RoleName][UNI,TOT]
$ws.Range("A19").Value = "[Roles]"
$ws.Range("B19").Value = "name"

$ws.Range("A20").Value = "Role"
$ws.Range("B20").Value = "RoleName"

$ws.Range("A21").Value = "Tutor"
$ws.Range("B21").Value = "Tutor"

$ws.Range("A22").Value = "Student"
$ws.Range("B22").Value = "Student"

$ws.Range("A23").Value = "GradStudent"
$ws.Range("B23").Value = "GradStudent"

$ws.Range("A24").Value = "AccountMgr"
$ws.Range("B24").Value = "AccountMgr"

# --- Selection / view tidy-up ------------------------------------------------
$ws.Range("B25").Select() | Out-Null
